$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Cell 1 (Tech Stack / left column of the skills table)
# ---------------------------------------------------------------------------

# 1) "AWS (EKS, EC2, Lambdas and other)" -> "AWS (CloudFormation, ECS, EKS, EC2, Lambdas and other)"
#    Plain text swap - keeps the existing bold run/formatting untouched.
$d.Content.Find.Execute( `
    "AWS (EKS, EC2, Lambdas and other)", $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "AWS (CloudFormation, ECS, EKS, EC2, Lambdas and other)", 2)

# 2) Drop " / Google Cloud Platform" (which used to trail the AWS line) and
#    replace it with a manual line break followed by "ElasticSearch / OpenSearch".
#    "^l" is Word's manual-line-break replacement code, so this yields a real
#    <w:br/> just like the ones already used throughout the document.
$d.Content.Find.Execute( `
    " / Google Cloud Platform", $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "^lElasticSearch / OpenSearch", 2)

# 2b) The whole inserted chunk above inherited the (non-bold) formatting of the
#     leading " / " it replaced, so re-bold just the two product names.
$rngElastic = $d.Content
$rngElastic.Find.Execute("ElasticSearch") | Out-Null
$rngElastic.Font.Bold = 1

$rngOpen = $d.Content
$rngOpen.Find.Execute("OpenSearch") | Out-Null
$rngOpen.Font.Bold = 1

# 2c) The paragraph used to end with a trailing line break after
#     "Google Cloud Platform"; the new last item ("OpenSearch") should not be
#     followed by one, so remove it.
$rngTail = $d.Content
$rngTail.Find.Execute("OpenSearch") | Out-Null
$rngTail.Collapse(0)
$rngTail.MoveEnd(1, 1)
$rngTail.Delete()

# ---------------------------------------------------------------------------
# Cell 2 (right column of the skills table)
# ---------------------------------------------------------------------------

# 3) Insert "Google Cloud Platform " (bold, with a plain trailing space) right
#    before "some experience with Kafka".
$d.Content.Find.Execute( `
    "some experience with Kafka", $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "Google Cloud Platform some experience with Kafka", 2)

$rngGcp = $d.Content
$rngGcp.Find.Execute("Google Cloud Platform some experience with Kafka") | Out-Null
$gcpSpaceStart = $rngGcp.Start + [string]"Google Cloud Platform".Length
$rngGcpSpace = $d.Range($gcpSpaceStart, $gcpSpaceStart + 1)
$rngGcpSpace.Font.Bold = 0

# 4) Append a trailing manual line break after "prefer functional programming
#    approach" (the paragraph used to end right after that text).
$d.Content.Find.Execute( `
    "prefer functional programming approach", $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "prefer functional programming approach^l", 2)

Write-Output "done"
